# Investigacion de reportes con graficos
# Adds a follow-up note (with a link to an ireport charts tutorial) to the
# "producto comodin" task row on Hoja1, and moves the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 42 ("producto comodin") gains a responsible person, a 100% status,
# and a note column (D) pointing at a reference article.
$ws.Range("B42").Value = "Agustina"
$ws.Range("C42").Value = 1
$ws.Range("C42").NumberFormat = "0%"
$ws.Range("D42").Value = "http://mygnet.net/articulos/java/creacion_de_graficos_en_ireport.707"

# Move the active selection to D43, matching where the author left off.
$ws.Activate()
$ws.Range("D43").Select()
